# Insert two new weekly price rows for "Betarraga" (row 750 and 751), pushing
# the existing rows 750-805 down to 752-807 (dimension grows from R805 to R807).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 750.. down by two to make room for the new entries.
$ws.Rows("750:751").Insert()

# New row 750
$ws.Range("A750").Value = 7
$ws.Range("B750").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C750").Value = "Ñuble"
$ws.Range("D750").Value = 45223
$ws.Range("E750").Value = 16
$ws.Range("F750").Value = 100114014
$ws.Range("G750").Value = "Betarraga"
$ws.Range("H750").Value = "Sin especificar"
$ws.Range("I750").Value = "Primera"
$ws.Range("J750").Value = 500
$ws.Range("K750").Value = 900
$ws.Range("L750").Value = 900
$ws.Range("M750").Value = 900
$ws.Range("N750").Value = "`$/paquete 5 unidades"
$ws.Range("O750").Value = "Provincia de Diguillín"
$ws.Range("P750").Value = 180
$ws.Range("Q750").Value = 5
$ws.Range("R750").Value = "Hortaliza"

# New row 751
$ws.Range("A751").Value = 7
$ws.Range("B751").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C751").Value = "Ñuble"
$ws.Range("D751").Value = 45223
$ws.Range("E751").Value = 16
$ws.Range("F751").Value = 100114014
$ws.Range("G751").Value = "Betarraga"
$ws.Range("H751").Value = "Sin especificar"
$ws.Range("I751").Value = "Segunda"
$ws.Range("J751").Value = 500
$ws.Range("K751").Value = 700
$ws.Range("L751").Value = 700
$ws.Range("M751").Value = 700
$ws.Range("N751").Value = "`$/paquete 5 unidades"
$ws.Range("O751").Value = "Provincia de Diguillín"
$ws.Range("P751").Value = 140
$ws.Range("Q751").Value = 5
$ws.Range("R751").Value = "Hortaliza"
